$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...chrome devtools. " -> "...chrome dev-tools. "
#    Replace the whole phrase (spanning the old proofErr-wrapped "devtools"
#    run) so the stale spell-check proofErr markers around "devtools" are
#    dropped, then re-split the run at the hyphen by toggling a character
#    property so the save-out produces three runs: "...chrome dev",
#    "-", "tools. " (mirroring how Word splits a run where a character was
#    typed in the middle of a word).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("chrome devtools. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "chrome dev-tools. ", 2)

$r = $d.Content
$r.Find.Execute("dev-tools") | Out-Null
$dashStart = $r.Start + 3
$rDash = $d.Range($dashStart, $dashStart + 1)
$rDash.Font.Bold = $true
$rDash.Font.Bold = $false

# ---------------------------------------------------------------------------
# 2) "extention" -> "extension"
#    Replace across both proofErr boundaries ("an extention" .. "name") so
#    the spellStart/spellEnd markers around "extention" are dropped, then
#    re-split the run around "extension" so the untouched neighboring text
#    ("c) we can use an " / " name \u201cReact Developer Tools\u201d") keeps
#    living in its own runs, just like before the edit.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("an extention name", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "an extension name", 2)

$r2 = $d.Content
$r2.Find.Execute("extension") | Out-Null
$rExt = $d.Range($r2.Start, $r2.End)
$rExt.Font.Bold = $true
$rExt.Font.Bold = $false
